$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text columns (B, C): the InternetComputer(DFINITY)/Maker rows swap rank ---
$ws.Range('B33').Value = 'Maker'
$ws.Range('C33').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'

# --- Numeric-looking text columns (D, E): force each cell to Text format so
#     values like "211.59" / "0.0500" are stored verbatim (not reinterpreted as
#     numbers). NumberFormat must be set per-cell (a multi-area union Range only
#     applies NumberFormat to its first area). Style is restored to "Normal"
#     afterwards so no stray style index is left on the cell. ---
$deAddrs = @('D2', 'E2', 'D3', 'E3', 'E4', 'D5', 'E6', 'E7', 'E8', 'D9', 'E9', 'E10', 'E11', 'D12', 'E12', 'D13', 'E13', 'E14', 'E15', 'E16', 'D17', 'E17', 'E18', 'E19', 'D20', 'E20', 'E21', 'E22', 'E23', 'D24', 'E24', 'D25', 'E25', 'E26', 'D27', 'E27', 'E29', 'D30', 'E30', 'E31', 'D32', 'E32', 'D33', 'E33', 'D34', 'E34', 'E35', 'E36', 'D37', 'E37', 'D38', 'E38', 'E39', 'E40', 'E41', 'D42', 'E42', 'E43', 'E44', 'D45', 'E45', 'D46', 'E46', 'D47', 'E47', 'D48', 'E48', 'D49', 'E49', 'D50', 'E50', 'D51', 'E51')
foreach ($addr in $deAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.185.73'
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').Value = '1.586.38'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '211.59'
$ws.Range('E6').Value = '  +0.74%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +0.37%  '
$ws.Range('D9').Value = '0.0606'
$ws.Range('E9').Value = '  -0.74%  '
$ws.Range('E10').Value = '  -1.49%  '
$ws.Range('E11').Value = '  +0.53%  '
$ws.Range('D12').Value = '1.809.47'
$ws.Range('E12').Value = '  +0.45%  '
$ws.Range('D13').Value = '1.598.13'
$ws.Range('E13').Value = '  +0.94%  '
$ws.Range('E14').Value = '  -1.15%  '
$ws.Range('E15').Value = '  +0.25%  '
$ws.Range('E16').Value = '  -0.32%  '
$ws.Range('D17').Value = '26.192.66'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('E18').Value = '  -0.39%  '
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('D20').Value = '213.10'
$ws.Range('E20').Value = '  +1.42%  '
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('E22').Value = '  -0.30%  '
$ws.Range('E23').Value = '  +0.34%  '
$ws.Range('D24').Value = '8.98'
$ws.Range('E24').Value = '  +1.72%  '
$ws.Range('D25').Value = '143.50'
$ws.Range('E25').Value = '  -0.47%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').Value = '6.99'
$ws.Range('E27').Value = '  -0.35%  '
$ws.Range('E29').Value = '  -0.72%  '
$ws.Range('D30').Value = '0.0496'
$ws.Range('E30').Value = '  -1.95%  '
$ws.Range('E31').Value = '  +0.70%  '
$ws.Range('D32').Value = '3.20'
$ws.Range('E32').Value = '  -0.90%  '
$ws.Range('D33').Value = '1.343.36'
$ws.Range('E33').Value = '  +4.61%  '
$ws.Range('D34').Value = '2.94'
$ws.Range('E34').Value = '  -1.84%  '
$ws.Range('E35').Value = '  +0.09%  '
$ws.Range('E36').Value = '  -1.16%  '
$ws.Range('D37').Value = '0.580'
$ws.Range('E37').Value = '  -3.90%  '
$ws.Range('D38').Value = '0.0167'
$ws.Range('E38').Value = '  +0.08%  '
$ws.Range('E39').Value = '  +0.88%  '
$ws.Range('E40').Value = '  +3.65%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('D42').Value = '0.951'
$ws.Range('E42').Value = '  -15.20%  '
$ws.Range('E43').Value = '  +0.74%  '
$ws.Range('E44').Value = '  +0.29%  '
$ws.Range('D45').Value = '1.721.95'
$ws.Range('E45').Value = '  +0.45%  '
$ws.Range('D46').Value = '61.16'
$ws.Range('E46').Value = '  -1.94%  '
$ws.Range('D47').Value = '86.08'
$ws.Range('E47').Value = '  -2.79%  '
$ws.Range('D48').Value = '0.0₆0102'
$ws.Range('E48').Value = '  +0.23%  '
$ws.Range('D49').Value = '1.48'
$ws.Range('E49').Value = '  -1.62%  '
$ws.Range('D50').Value = '0.0987'
$ws.Range('E50').Value = '  -1.88%  '
$ws.Range('D51').Value = '0.0500'
$ws.Range('E51').Value = '  -0.90%  '

foreach ($addr in $deAddrs) {
    $ws.Range($addr).Style = "Normal"
}
